# Update the header text for the gene-name column.
# "Name of the gene" -> "Gene" (reflects that stretches can now span
# multiple genes, per the commit message).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E1").Value = "Gene"
